$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at position 2 (shifts existing rows 2-16 down to 3-17)
# ---------------------------------------------------------------------------
$ws.Rows(2).Insert()

# Populate the new row 2 with the newest filing (10-Q, filed 2024-08-13,
# reporting period 2024-06-30).
$ws.Range("A2").Value = "10-Q"
$ws.Range("B2").Value = "Quarterly report [Sections 13 or 15(d)]"
$ws.Range("C2").Value = 45517
$ws.Range("D2").Value = 45473
$ws.Range("E2").Value = "https://www.sec.gov/Archives/edgar/data/1803498/000180349824000044/bcred-20240630.htm"

# Copy the formatting of the (now shifted) row that used to be row 2, cell by
# cell, so the new row looks the same as its neighbours (date format on
# C/D, hyperlink-like style on E, etc).
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("D3").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("E3").Copy()
$ws.Range("E2").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Update the html_link column (E) for the rows whose URL moved from the
#    old SEC "-index.htm" filing-index page to the actual filing document.
#    (Rows are now offset by +1 after the insert above.)
# ---------------------------------------------------------------------------
$ws.Range("E7").Value  = "https://www.sec.gov/Archives/edgar/data/1803498/000180349823000012/bcred-20230331.htm"
$ws.Range("E8").Value  = "https://www.sec.gov/Archives/edgar/data/1803498/000180349823000008/bcred-20221231.htm"
$ws.Range("E9").Value  = "https://www.sec.gov/Archives/edgar/data/1803498/000180349822000021/bcred-20220930.htm"
$ws.Range("E10").Value = "https://www.sec.gov/Archives/edgar/data/1803498/000180349822000017/bcred-20220630.htm"
$ws.Range("E11").Value = "https://www.sec.gov/Archives/edgar/data/1803498/000180349822000013/bcred-03312022x10q.htm"
$ws.Range("E12").Value = "https://www.sec.gov/Archives/edgar/data/1803498/000180349822000009/bcred-20211231x10k.htm"
$ws.Range("E13").Value = "https://www.sec.gov/Archives/edgar/data/1803498/000180349821000012/bcred-09302021x10q.htm"
$ws.Range("E14").Value = "https://www.sec.gov/Archives/edgar/data/1803498/000180349821000009/bcred-06302021x10q.htm"
$ws.Range("E15").Value = "https://www.sec.gov/Archives/edgar/data/1803498/000180349821000003/bcred-03312021x10q.htm"
$ws.Range("E16").Value = "https://www.sec.gov/Archives/edgar/data/1803498/000119312521069991/d145785d10k.htm"
$ws.Range("E17").Value = "https://www.sec.gov/Archives/edgar/data/1803498/000119312520294937/d50632d10q.htm"

# ---------------------------------------------------------------------------
# 3. Clear the two stray leftover header labels in I1/J1 (their underlying
#    data columns were removed long ago, only the header text remained).
# ---------------------------------------------------------------------------
$ws.Range("I1").ClearContents()
$ws.Range("J1").ClearContents()

# ---------------------------------------------------------------------------
# 4. Hyperlinks: only 3 live hyperlinks remain (the row that used to carry a
#    4th hyperlink lost it). Row-insert does not auto-shift the hyperlink
#    collection in this engine, so rebuild it explicitly in final position.
# ---------------------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E3"), "https://www.sec.gov/Archives/edgar/data/1803498/000180349824000022/bcred-20240331.htm")
$ws.Hyperlinks.Add($ws.Range("E4"), "https://www.sec.gov/Archives/edgar/data/1803498/000180349824000007/bcred-20231231.htm")
$ws.Hyperlinks.Add($ws.Range("E5"), "https://www.sec.gov/Archives/edgar/data/1803498/000180349823000021/bcred-20230930.htm")

# ---------------------------------------------------------------------------
# 5. Selection / active cell moves to E18 (mirrors the saved view state).
# ---------------------------------------------------------------------------
$ws.Range("E18").Select()
